$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 697.125
$ws.Range("I18").Value = 799.1667
$ws.Range("J18").Value = 391
$ws.Range("K18").Value = 799.1667
$ws.Range("L18").Value = 391
$ws.Range("M18").Value = -515.1667
$ws.Range("N18").Value = -959
$ws.Range("H97").Value = 10000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -30992
$ws.Range("H130").Value = 95500
$ws.Range("J130").Value = 95500
$ws.Range("L130").Value = 95500
$ws.Range("N130").Value = -105540
$ws.Range("H131").Value = 2135.8
$ws.Range("I131").Value = 1973.1111
$ws.Range("K131").Value = 5919.3333
$ws.Range("M131").Value = -879.3333000000002
$ws.Range("H132").Value = 2413.2444
$ws.Range("I132").Value = 2483.5
$ws.Range("K132").Value = 7450.5
$ws.Range("M132").Value = -4920.5
$ws.Range("H138").Value = 2133.1396
$ws.Range("I138").Value = 1540
$ws.Range("K138").Value = 4620
$ws.Range("M138").Value = 520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1713
$ws.Range("I2").Value = 1713
$ws.Range("K2").Value = 1713
$ws.Range("M2").Value = -1600
$ws.Range("H32").Value = 7468441.5
$ws.Range("I32").Value = 7817165
$ws.Range("K32").Value = 7817165
$ws.Range("M32").Value = -7816878
$ws.Range("H45").Value = 33334974
$ws.Range("I45").Value = 41667800
$ws.Range("K45").Value = 41667800
$ws.Range("M45").Value = -41667423
$ws.Range("H61").Value = 17281120
$ws.Range("I61").Value = 20004240
$ws.Range("J61").Value = 261614.5
$ws.Range("K61").Value = 20004240
$ws.Range("L61").Value = 261614.5
$ws.Range("M61").Value = -20004028
$ws.Range("N61").Value = -262038.5
$ws.Range("H63").Value = 6596
$ws.Range("I63").Value = 6422
$ws.Range("K63").Value = 6422
$ws.Range("M63").Value = -5736
$ws.Range("H66").Value = 6596
$ws.Range("I66").Value = 6422
$ws.Range("K66").Value = 32110
$ws.Range("M66").Value = -28678
$ws.Range("H74").Value = 6671998.5
$ws.Range("I74").Value = 9618367
$ws.Range("K74").Value = 9618367
$ws.Range("M74").Value = -9617493
$ws.Range("H77").Value = 6671998.5
$ws.Range("I77").Value = 9618367
$ws.Range("K77").Value = 48091835
$ws.Range("M77").Value = -48087467
$ws.Range("H80").Value = 48705.25
$ws.Range("I80").Value = 35000
$ws.Range("J80").Value = 53273.668
$ws.Range("K80").Value = 35000
$ws.Range("L80").Value = 53273.668
$ws.Range("M80").Value = -34002
$ws.Range("N80").Value = -55269.668
$ws.Range("H83").Value = 48705.25
$ws.Range("I83").Value = 35000
$ws.Range("J83").Value = 53273.668
$ws.Range("K83").Value = 105000
$ws.Range("L83").Value = 159821.004
$ws.Range("M83").Value = -100008
$ws.Range("N83").Value = -169805.004
$ws.Range("H110").Value = 1816
$ws.Range("I110").Value = 1816
$ws.Range("K110").Value = 1816
$ws.Range("M110").Value = 229
$ws.Range("H116").Value = 1713
$ws.Range("I116").Value = 1713
$ws.Range("K116").Value = 1713
$ws.Range("M116").Value = 581
$ws.Range("H122").Value = 4169
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4169
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12507
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17407
$ws.Range("H136").Value = 17281120
$ws.Range("I136").Value = 20004240
$ws.Range("J136").Value = 261614.5
$ws.Range("K136").Value = 60012720
$ws.Range("L136").Value = 784843.5
$ws.Range("M136").Value = -60010170
$ws.Range("N136").Value = -789943.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1713
$ws.Range("I3").Value = 1713
$ws.Range("K3").Value = 1713
$ws.Range("M3").Value = -1599
$ws.Range("H82").Value = 37367.2
$ws.Range("I82").Value = 34677.8
$ws.Range("J82").Value = 40056.6
$ws.Range("K82").Value = 34677.8
$ws.Range("L82").Value = 40056.6
$ws.Range("M82").Value = -34294.8
$ws.Range("N82").Value = -40822.6
$ws.Range("H85").Value = 37367.2
$ws.Range("I85").Value = 34677.8
$ws.Range("J85").Value = 40056.6
$ws.Range("K85").Value = 34677.8
$ws.Range("L85").Value = 40056.6
$ws.Range("M85").Value = -33351.8
$ws.Range("N85").Value = -42708.6
$ws.Range("H99").Value = 5882.2
$ws.Range("I99").Value = 8484.933999999999
$ws.Range("J99").Value = 3279.4666
$ws.Range("K99").Value = 8484.933999999999
$ws.Range("L99").Value = 3279.4666
$ws.Range("M99").Value = -6986.933999999999
$ws.Range("N99").Value = -6275.4666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000.125
$ws.Range("J4").Value = 1000.125
$ws.Range("L4").Value = 1000.125
$ws.Range("N4").Value = -1224.125
$ws.Range("H31").Value = 789471.4
$ws.Range("I31").Value = 17341.615
$ws.Range("K31").Value = 17341.615
$ws.Range("M31").Value = -17046.615
$ws.Range("H34").Value = 789471.4
$ws.Range("I34").Value = 17341.615
$ws.Range("K34").Value = 17341.615
$ws.Range("M34").Value = -17139.615
$ws.Range("H58").Value = 1745.425
$ws.Range("I58").Value = 1354.0857
$ws.Range("J58").Value = 4484.8
$ws.Range("K58").Value = 1354.0857
$ws.Range("L58").Value = 4484.8
$ws.Range("M58").Value = -1151.0857
$ws.Range("N58").Value = -4890.8
$ws.Range("H112").Value = 85662
$ws.Range("J112").Value = 85662
$ws.Range("L112").Value = 85662
$ws.Range("N112").Value = -88616
$ws.Range("H132").Value = 2696
$ws.Range("I132").Value = 2613.0908
$ws.Range("K132").Value = 7839.2724
$ws.Range("M132").Value = -5309.2724
$ws.Range("H134").Value = 403894.3
$ws.Range("I134").Value = 557961.75
$ws.Range("K134").Value = 1673885.25
$ws.Range("M134").Value = -1671350.25
$ws.Range("H136").Value = 1745.425
$ws.Range("I136").Value = 1354.0857
$ws.Range("J136").Value = 4484.8
$ws.Range("K136").Value = 4062.2571
$ws.Range("L136").Value = 13454.4
$ws.Range("M136").Value = -1512.2571
$ws.Range("N136").Value = -18554.4
$ws.Range("H137").Value = 69999.336
$ws.Range("J137").Value = 69999.336
$ws.Range("L137").Value = 69999.336
$ws.Range("N137").Value = -80199.336
$ws.Range("H138").Value = 95999.5
$ws.Range("J138").Value = 95999.5
$ws.Range("L138").Value = 95999.5
$ws.Range("N138").Value = -106279.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 135.27777
$ws.Range("I2").Value = 38.22222
$ws.Range("J2").Value = 232.33333
$ws.Range("K2").Value = 38.22222
$ws.Range("L2").Value = 232.33333
$ws.Range("M2").Value = 74.77778000000001
$ws.Range("N2").Value = -458.33333
$ws.Range("H127").Value = 100140
$ws.Range("J127").Value = 100140
$ws.Range("L127").Value = 100140
$ws.Range("N127").Value = -110060
$ws.Range("H132").Value = 43487744
$ws.Range("I132").Value = 58827490
$ws.Range("J132").Value = 25135.5
$ws.Range("K132").Value = 176482470
$ws.Range("L132").Value = 75406.5
$ws.Range("M132").Value = -176479940
$ws.Range("N132").Value = -80466.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 749.375
$ws.Range("I22").Value = 899.4
$ws.Range("J22").Value = 499.33334
$ws.Range("K22").Value = 899.4
$ws.Range("L22").Value = 499.33334
$ws.Range("M22").Value = -604.4
$ws.Range("N22").Value = -1089.33334
$ws.Range("H27").Value = 749.375
$ws.Range("I27").Value = 899.4
$ws.Range("J27").Value = 499.33334
$ws.Range("K27").Value = 899.4
$ws.Range("L27").Value = 499.33334
$ws.Range("M27").Value = -792.4
$ws.Range("N27").Value = -713.33334
$ws.Range("H136").Value = 54963.88
$ws.Range("I136").Value = 5522.7144
$ws.Range("J136").Value = 314530
$ws.Range("K136").Value = 16568.1432
$ws.Range("L136").Value = 943590
$ws.Range("M136").Value = -14018.1432
$ws.Range("N136").Value = -948690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5253750
$ws.Range("I5").Value = 500000
$ws.Range("K5").Value = 500000
$ws.Range("M5").Value = -499888
$ws.Range("H14").Value = 3332.5557
$ws.Range("H33").Value = 20000
$ws.Range("J33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("N33").Value = -20500
$ws.Range("H36").Value = 20000
$ws.Range("J36").Value = 20000
$ws.Range("L36").Value = 20000
$ws.Range("N36").Value = -20500
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H40").Value = 30495
$ws.Range("J40").Value = 30495
$ws.Range("L40").Value = 30495
$ws.Range("N40").Value = -30793
$ws.Range("H122").Value = 1451.1666
$ws.Range("I122").Value = 1401.2727
$ws.Range("K122").Value = 4203.8181
$ws.Range("M122").Value = -1753.8181
$ws.Range("H132").Value = 5169.15
$ws.Range("I132").Value = 4780.9375
$ws.Range("J132").Value = 6722
$ws.Range("K132").Value = 14342.8125
$ws.Range("L132").Value = 20166
$ws.Range("M132").Value = -11812.8125
$ws.Range("N132").Value = -25226
